# Add a new "Croatia" Test Data sheet, cloned from the existing "Turkey"
# sheet (same layout/styles), with the Croatia-specific values filled in,
# placed as the last (rightmost) tab and made the active sheet.

$wb = $excel.ActiveWorkbook

$turkey = $wb.Worksheets.Item("Turkey")

# Duplicate "Turkey" immediately after itself -> becomes the new last sheet
# and the active sheet (mirrors Excel's own "Move or Copy... (Create a
# copy)" behaviour).
$turkey.Copy($null, $turkey)
$croatia = $wb.ActiveSheet
$croatia.Name = "Croatia"

# Fill in the market-specific cells for Croatia.
$croatia.Range("B2").Value = "Croatia Market"
$croatia.Range("B4").Value = "NGC-3139/T2421"

# Turkey is no longer the selected tab; restore its cursor to a whole-sheet
# selection (matches the target state) before moving focus to Croatia.
$turkey.Range("A1:XFD1048576").Select()

$croatia.Activate()
$croatia.Range("B4").Select()
